$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Alignment constants (xlCenter = -4108)
$xlCenter = -4108

# --- Row 6 & 7: História placeholder rows (only column A filled) ---
$ws.Range("A6").Value = "7 - Permitir o usuário avaliar com gostei/não gostei"
$ws.Range("A6").HorizontalAlignment = $xlCenter
$ws.Range("A6").VerticalAlignment = $xlCenter

$ws.Range("A7").Value = "7 - Permitir o usuário avaliar com gostei/não gostei"
$ws.Range("A7").HorizontalAlignment = $xlCenter
$ws.Range("A7").VerticalAlignment = $xlCenter

# --- Column A for rows 8,9,10 (story 13) and 11 (story 8) ---
$ws.Range("A8").Value = "13 - Permitir o usuário favoritar cursos"
$ws.Range("A8").HorizontalAlignment = $xlCenter
$ws.Range("A8").VerticalAlignment = $xlCenter

$ws.Range("A9").Value = "13 - Permitir o usuário favoritar cursos"
$ws.Range("A9").HorizontalAlignment = $xlCenter
$ws.Range("A9").VerticalAlignment = $xlCenter

$ws.Range("A10").Value = "13 - Permitir o usuário favoritar cursos"
$ws.Range("A10").HorizontalAlignment = $xlCenter
$ws.Range("A10").VerticalAlignment = $xlCenter

$ws.Range("A11").Value = "8 - Permitir o usuário avaliar um curso com comentário"
$ws.Range("A11").HorizontalAlignment = $xlCenter
$ws.Range("A11").VerticalAlignment = $xlCenter

# --- Row 11 (C then B, matching author's original entry order) ---
$ws.Range("C11").Value = "Clicando no curso após a pesquisa, foi mostrado a tela de detalhes do curso, na caixa de texto disponível, foi digitado o comentário e ao postar, a tela foi atualizada mostrando o comentário e todos os outros que o curso possuia"
$ws.Range("C11").HorizontalAlignment = $xlCenter
$ws.Range("C11").VerticalAlignment = $xlCenter
$ws.Range("C11").WrapText = $true

$ws.Range("B11").Value = "Ao clicar em um curso, depois da pesquisa ser efetuada, o usuário é redirecionado para a página dos detalhes do curso, onde estará disponível para ele uma caixa de texto e um botão para submeter o comentário, após a realização do mesmo, a tela é atualizada com o comentário recém postado e os demais"
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").HorizontalAlignment = $xlCenter
$ws.Range("B11").VerticalAlignment = $xlCenter
$ws.Range("B11").WrapText = $true

# --- Row 8 & 9 (B then C) ---
$ws.Range("B8").Value = "Ao clicar em um curso, depois da pesquisa ser efetuada, o usuário é redirecionado para a página dos detalhes do curso, tela esta que deve disponibilizar o botão para o usuário favoritar o curso que está na tela"
$ws.Range("B8").HorizontalAlignment = $xlCenter
$ws.Range("B8").VerticalAlignment = $xlCenter
$ws.Range("B8").WrapText = $true

$ws.Range("B9").Value = "Ao clicar em um curso, depois da pesquisa ser efetuada, o usuário é redirecionado para a página dos detalhes do curso, tela esta que deve disponibilizar o botão para o usuário favoritar o curso que está na tela"
$ws.Range("B9").HorizontalAlignment = $xlCenter
$ws.Range("B9").VerticalAlignment = $xlCenter
$ws.Range("B9").WrapText = $true

$ws.Range("C8").Value = "Clicando no curso após a pesquisa, foi mostrado a tela de detalhes do curso, e também o botão para favoritar o curso"
$ws.Range("C8").HorizontalAlignment = $xlCenter
$ws.Range("C8").VerticalAlignment = $xlCenter
$ws.Range("C8").WrapText = $true

$ws.Range("C9").Value = "Clicando no curso após a pesquisa, foi mostrado a tela de detalhes do curso, e também o botão para favoritar o curso"
$ws.Range("C9").HorizontalAlignment = $xlCenter
$ws.Range("C9").VerticalAlignment = $xlCenter
$ws.Range("C9").WrapText = $true

# --- Row 10 (B then C) ---
$ws.Range("B10").Value = "Ao clicar na aba de cursos favoritados o aplicativo deve listar todos os cursos que o usuário logado favoritou e se não tiver usuário logado, sugerir para que faça o login"
$ws.Range("B10").HorizontalAlignment = $xlCenter
$ws.Range("B10").VerticalAlignment = $xlCenter
$ws.Range("B10").WrapText = $true

$ws.Range("C10").Value = "Indo na aba de cursos favoritados, foi listado todos que o usuário logado favoritou e sem login foi sugerido para que o usuário o faça"
$ws.Range("C10").HorizontalAlignment = $xlCenter
$ws.Range("C10").VerticalAlignment = $xlCenter
$ws.Range("C10").WrapText = $true

# --- Status column (D) ---
$ws.Range("D8").Value = "Aprovado"
$ws.Range("D8").Style = "Bom"
$ws.Range("D8").HorizontalAlignment = $xlCenter
$ws.Range("D8").VerticalAlignment = $xlCenter

$ws.Range("D9").Value = "Aprovado"
$ws.Range("D9").Style = "Bom"
$ws.Range("D9").HorizontalAlignment = $xlCenter
$ws.Range("D9").VerticalAlignment = $xlCenter

$ws.Range("D10").Value = "Pendente"
$ws.Range("D10").Style = "Neutro"
$ws.Range("D10").HorizontalAlignment = $xlCenter
$ws.Range("D10").VerticalAlignment = $xlCenter

$ws.Range("D11").Value = "Pendente"
$ws.Range("D11").Style = "Neutro"
$ws.Range("D11").HorizontalAlignment = $xlCenter
$ws.Range("D11").VerticalAlignment = $xlCenter

# --- Row heights ---
$ws.Rows.Item(8).RowHeight = 45
$ws.Rows.Item(9).RowHeight = 45
$ws.Rows.Item(10).RowHeight = 45
$ws.Rows.Item(11).RowHeight = 60

# --- Column A width ---
$ws.Columns.Item(1).ColumnWidth = 58.42578125

# --- Selection change ---
$ws.Range("F5").Select()
